# Update the PSSM score matrix (B2:K21) with the recomputed values from the
# supplemental-figures rerun. Row/column headers (row 1, column A) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 20,10
$data[0,0] = -19.0161961917914
$data[0,1] = 2.539829397197745
$data[0,2] = -19.0161961917914
$data[0,3] = -19.0161961917914
$data[0,4] = -19.0161961917914
$data[0,5] = -19.0161961917914
$data[0,6] = -19.0161961917914
$data[0,7] = -19.0161961917914
$data[0,8] = -19.0161961917914
$data[0,9] = -19.0161961917914
$data[1,0] = -19.0161961917914
$data[1,1] = -19.0161961917914
$data[1,2] = -19.0161961917914
$data[1,3] = -19.0161961917914
$data[1,4] = -19.0161961917914
$data[1,5] = -19.0161961917914
$data[1,6] = -19.0161961917914
$data[1,7] = -19.0161961917914
$data[1,8] = -19.0161961917914
$data[1,9] = -19.0161961917914
$data[2,0] = -19.0161961917914
$data[2,1] = 2.240402392136413
$data[2,2] = 2.936184259753373
$data[2,3] = -19.0161961917914
$data[2,4] = 2.570138307204047
$data[2,5] = -19.0161961917914
$data[2,6] = 1.862358446750684
$data[2,7] = -19.0161961917914
$data[2,8] = 2.257791152835225
$data[2,9] = -19.0161961917914
$data[3,0] = -19.0161961917914
$data[3,1] = 0.9215353725044028
$data[3,2] = -19.0161961917914
$data[3,3] = -19.0161961917914
$data[3,4] = -19.0161961917914
$data[3,5] = 2.204038642474678
$data[3,6] = -19.0161961917914
$data[3,7] = -19.0161961917914
$data[3,8] = -19.0161961917914
$data[3,9] = -19.0161961917914
$data[4,0] = -19.0161961917914
$data[4,1] = -19.0161961917914
$data[4,2] = -19.0161961917914
$data[4,3] = -19.0161961917914
$data[4,4] = -19.0161961917914
$data[4,5] = -19.0161961917914
$data[4,6] = -19.0161961917914
$data[4,7] = -19.0161961917914
$data[4,8] = -19.0161961917914
$data[4,9] = -19.0161961917914
$data[5,0] = 2.957330571320866
$data[5,1] = -19.0161961917914
$data[5,2] = -19.0161961917914
$data[5,3] = -19.0161961917914
$data[5,4] = -19.0161961917914
$data[5,5] = -19.0161961917914
$data[5,6] = -19.0161961917914
$data[5,7] = -19.0161961917914
$data[5,8] = -19.0161961917914
$data[5,9] = -19.0161961917914
$data[6,0] = -19.0161961917914
$data[6,1] = -19.0161961917914
$data[6,2] = -19.0161961917914
$data[6,3] = 2.980014185121259
$data[6,4] = -19.0161961917914
$data[6,5] = -19.0161961917914
$data[6,6] = -19.0161961917914
$data[6,7] = -19.0161961917914
$data[6,8] = -19.0161961917914
$data[6,9] = -19.0161961917914
$data[7,0] = 3.612719418219589
$data[7,1] = -19.0161961917914
$data[7,2] = -19.0161961917914
$data[7,3] = -19.0161961917914
$data[7,4] = -19.0161961917914
$data[7,5] = -19.0161961917914
$data[7,6] = -19.0161961917914
$data[7,7] = -19.0161961917914
$data[7,8] = -19.0161961917914
$data[7,9] = -19.0161961917914
$data[8,0] = -19.0161961917914
$data[8,1] = -19.0161961917914
$data[8,2] = -19.0161961917914
$data[8,3] = -19.0161961917914
$data[8,4] = -19.0161961917914
$data[8,5] = -19.0161961917914
$data[8,6] = -19.0161961917914
$data[8,7] = -19.0161961917914
$data[8,8] = -19.0161961917914
$data[8,9] = 2.234586231177514
$data[9,0] = -19.0161961917914
$data[9,1] = -19.0161961917914
$data[9,2] = -19.0161961917914
$data[9,3] = 2.009365837978404
$data[9,4] = -19.0161961917914
$data[9,5] = 2.825069073446648
$data[9,6] = -19.0161961917914
$data[9,7] = -19.0161961917914
$data[9,8] = -19.0161961917914
$data[9,9] = 1.661668803355541
$data[10,0] = -19.0161961917914
$data[10,1] = -19.0161961917914
$data[10,2] = -19.0161961917914
$data[10,3] = -19.0161961917914
$data[10,4] = -19.0161961917914
$data[10,5] = -19.0161961917914
$data[10,6] = -19.0161961917914
$data[10,7] = -19.0161961917914
$data[10,8] = -19.0161961917914
$data[10,9] = -19.0161961917914
$data[11,0] = -19.0161961917914
$data[11,1] = -19.0161961917914
$data[11,2] = -19.0161961917914
$data[11,3] = 1.620705399668374
$data[11,4] = -19.0161961917914
$data[11,5] = -19.0161961917914
$data[11,6] = -19.0161961917914
$data[11,7] = -19.0161961917914
$data[11,8] = 2.360018773724543
$data[11,9] = 1.664881759838748
$data[12,0] = -19.0161961917914
$data[12,1] = -19.0161961917914
$data[12,2] = 1.713614547465644
$data[12,3] = -19.0161961917914
$data[12,4] = -19.0161961917914
$data[12,5] = -19.0161961917914
$data[12,6] = -19.0161961917914
$data[12,7] = -19.0161961917914
$data[12,8] = -19.0161961917914
$data[12,9] = 2.10449541846273
$data[13,0] = -19.0161961917914
$data[13,1] = -19.0161961917914
$data[13,2] = -0.1644075070577134
$data[13,3] = -19.0161961917914
$data[13,4] = -19.0161961917914
$data[13,5] = -19.0161961917914
$data[13,6] = -19.0161961917914
$data[13,7] = -19.0161961917914
$data[13,8] = -19.0161961917914
$data[13,9] = -19.0161961917914
$data[14,0] = -19.0161961917914
$data[14,1] = -19.0161961917914
$data[14,2] = -19.0161961917914
$data[14,3] = -19.0161961917914
$data[14,4] = -19.0161961917914
$data[14,5] = -19.0161961917914
$data[14,6] = -19.0161961917914
$data[14,7] = -19.0161961917914
$data[14,8] = 2.165437140466991
$data[14,9] = -19.0161961917914
$data[15,0] = -19.0161961917914
$data[15,1] = 0.7970542229597187
$data[15,2] = 0.07184924949575146
$data[15,3] = -19.0161961917914
$data[15,4] = -19.0161961917914
$data[15,5] = -19.0161961917914
$data[15,6] = 1.036760509473488
$data[15,7] = -19.0161961917914
$data[15,8] = 1.479566142521006
$data[15,9] = -19.0161961917914
$data[16,0] = -19.0161961917914
$data[16,1] = -19.0161961917914
$data[16,2] = -19.0161961917914
$data[16,3] = -19.0161961917914
$data[16,4] = -19.0161961917914
$data[16,5] = -19.0161961917914
$data[16,6] = 0.8131874428622797
$data[16,7] = -19.0161961917914
$data[16,8] = 1.490027115466878
$data[16,9] = -19.0161961917914
$data[17,0] = -19.0161961917914
$data[17,1] = -19.0161961917914
$data[17,2] = 1.631477415660003
$data[17,3] = -19.0161961917914
$data[17,4] = -19.0161961917914
$data[17,5] = -19.0161961917914
$data[17,6] = 1.917780388866621
$data[17,7] = -19.0161961917914
$data[17,8] = -19.0161961917914
$data[17,9] = -19.0161961917914
$data[18,0] = -19.0161961917914
$data[18,1] = 1.485836280070299
$data[18,2] = 2.008865622991201
$data[18,3] = -19.0161961917914
$data[18,4] = 3.813658329659574
$data[18,5] = -19.0161961917914
$data[18,6] = 1.88909458294458
$data[18,7] = 4.321925509931832
$data[18,8] = -19.0161961917914
$data[18,9] = 2.219872321645882
$data[19,0] = -19.0161961917914
$data[19,1] = 1.597892640256087
$data[19,2] = -19.0161961917914
$data[19,3] = 2.324443373301531
$data[19,4] = -19.0161961917914
$data[19,5] = 3.054121655063636
$data[19,6] = 2.34278599935432
$data[19,7] = -19.0161961917914
$data[19,8] = -19.0161961917914
$data[19,9] = -19.0161961917914
$ws.Range("B2:K21").Value = $data
